$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

# Simple single-value replacements (rows are 1-indexed)
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "34"
Set-CellText $t 6 "0.00070"
Set-CellText $t 7 "0.00021"
Set-CellText $t 9 "0.00037"
Set-CellText $t 10 "0.00047"
Set-CellText $t 11 "0.00053"
Set-CellText $t 12 "0.00878"

# Rows 44-46 collapse multi-run tab-separated content into a single value
Set-CellText $t 44 "99.99"
Set-CellText $t 45 "0.01"
Set-CellText $t 46 "65"
